$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 was a blank placeholder row (all cells empty) - clear it so it
# disappears from the written sheetData entirely (rows 3+ keep their
# original row numbers).
$ws.Range("A2:L2").ClearContents()

# Row 4 ("Motivo Ingreso" / "Motivo Retiro" columns) no longer carries a
# (blank) value for this closing entry - remove those two cells.
$ws.Range("F4").ClearContents()
$ws.Range("H4").ClearContents()

# New closing record appended as row 5. The sheet stores every figure as
# literal text (e.g. "2000.00"), so force text formatting before writing
# the numeric-looking values, then drop back to the default style so no
# extra formatting is left behind on the cell.
$row5TextCells = "B5","C5","D5","E5","G5","I5"
foreach ($addr in $row5TextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("A5").Value = "2025-11-06 01:04:03"
$ws.Range("B5").Value = "2000.00"
$ws.Range("C5").Value = "1000.00"
$ws.Range("D5").Value = "29999.00"
$ws.Range("E5").Value = "28999.00"
$ws.Range("G5").Value = "0.00"
$ws.Range("I5").Value = "0.00"
foreach ($addr in $row5TextCells) {
    $ws.Range($addr).Style = "Normal"
}

# Another closing record appended as row 6.
$row6TextCells = "B6","C6","D6","E6","G6","I6"
foreach ($addr in $row6TextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("A6").Value = "2025-11-06 01:36:21"
$ws.Range("B6").Value = "4500.00"
$ws.Range("C6").Value = "1000.00"
$ws.Range("D6").Value = "9000.00"
$ws.Range("E6").Value = "5500.00"
$ws.Range("G6").Value = "0.00"
$ws.Range("I6").Value = "0.00"
foreach ($addr in $row6TextCells) {
    $ws.Range($addr).Style = "Normal"
}
